$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 233.66667
$ws.Range("I6").Value = 233.66667
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 701.00001
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -589.00001
$ws.Range("N6").ClearContents()
$ws.Range("H17").Value = 1143573
$ws.Range("I17").Value = 142.5
$ws.Range("J17").Value = 1175782.4
$ws.Range("K17").Value = 427.5
$ws.Range("L17").Value = 3527347.2
$ws.Range("M17").Value = -259.5
$ws.Range("N17").Value = -3527683.2
$ws.Range("H19").Value = 1274.3914
$ws.Range("I19").Value = 2498.2
$ws.Range("J19").Value = 934.44446
$ws.Range("K19").Value = 2498.2
$ws.Range("L19").Value = 934.44446
$ws.Range("M19").Value = -2323.2
$ws.Range("N19").Value = -1284.44446
$ws.Range("H64").Value = 3452.9443
$ws.Range("I64").Value = 3416.6667
$ws.Range("J64").Value = 3489.2222
$ws.Range("K64").Value = 3416.6667
$ws.Range("L64").Value = 3489.2222
$ws.Range("M64").Value = -3168.6667
$ws.Range("N64").Value = -3985.2222
$ws.Range("H67").Value = 3452.9443
$ws.Range("I67").Value = 3416.6667
$ws.Range("J67").Value = 3489.2222
$ws.Range("K67").Value = 3416.6667
$ws.Range("L67").Value = 3489.2222
$ws.Range("M67").Value = -2558.6667
$ws.Range("N67").Value = -5205.2222
$ws.Range("H96").Value = 949.125
$ws.Range("I96").Value = 966.46155
$ws.Range("J96").Value = 928.63635
$ws.Range("K96").Value = 2899.38465
$ws.Range("L96").Value = 2785.90905
$ws.Range("M96").Value = -1526.38465
$ws.Range("N96").Value = -5531.90905
$ws.Range("H139").Value = 38095.6
$ws.Range("J139").Value = 38095.6
$ws.Range("L139").Value = 38095.6
$ws.Range("N139").Value = -48375.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 9262132
$ws.Range("I132").Value = 17859750
$ws.Range("J132").Value = 3158
$ws.Range("K132").Value = 53579250
$ws.Range("L132").Value = 9474
$ws.Range("M132").Value = -53576720
$ws.Range("N132").Value = -14534

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 41593.668
$ws.Range("I75").Value = 15666.667
$ws.Range("J75").Value = 50236
$ws.Range("K75").Value = 15666.667
$ws.Range("L75").Value = 50236
$ws.Range("M75").Value = -14730.667
$ws.Range("N75").Value = -52108
$ws.Range("H78").Value = 41593.668
$ws.Range("I78").Value = 15666.667
$ws.Range("J78").Value = 50236
$ws.Range("K78").Value = 47000.001
$ws.Range("L78").Value = 150708
$ws.Range("M78").Value = -42320.001
$ws.Range("N78").Value = -160068
$ws.Range("H97").Value = 1990
$ws.Range("I97").Value = 1990
$ws.Range("K97").Value = 1990
$ws.Range("M97").Value = -999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 13891366
$ws.Range("I132").Value = 17859010
$ws.Range("K132").Value = 53577030
$ws.Range("M132").Value = -53574500
$ws.Range("H140").Value = 37510.57
$ws.Range("J140").Value = 37510.57
$ws.Range("L140").Value = 37510.57
$ws.Range("N140").Value = -47870.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 25000120
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H23").Value = 147.64285
$ws.Range("I23").Value = 100.25
$ws.Range("K23").Value = 300.75
$ws.Range("M23").Value = -65.75
$ws.Range("H34").Value = 1274.238
$ws.Range("J34").Value = 1521.7059
$ws.Range("L34").Value = 4565.1177
$ws.Range("N34").Value = -4733.1177
$ws.Range("H47").Value = 293.42856
$ws.Range("I47").Value = 231.2
$ws.Range("J47").Value = 449
$ws.Range("K47").Value = 693.5999999999999
$ws.Range("L47").Value = 1347
$ws.Range("M47").Value = -262.5999999999999
$ws.Range("N47").Value = -2209
$ws.Range("H80").Value = 2993.3333
$ws.Range("J80").Value = 2993.3333
$ws.Range("L80").Value = 8979.999899999999
$ws.Range("N80").Value = -10851.9999
$ws.Range("H83").Value = 2993.3333
$ws.Range("J83").Value = 2993.3333
$ws.Range("L83").Value = 26939.9997
$ws.Range("N83").Value = -36299.9997
$ws.Range("H118").Value = 1239.8948
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("H131").Value = 833.62
$ws.Range("J131").Value = 858.7311999999999
$ws.Range("L131").Value = 2576.1936
$ws.Range("N131").Value = -12656.1936

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 3000
$ws.Range("I35").Value = 3000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 3000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -2702
$ws.Range("N35").ClearContents()
$ws.Range("H70").Value = 300000
$ws.Range("I70").Value = 300000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 300000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -299730
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 300000
$ws.Range("I73").Value = 300000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 300000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -299064
$ws.Range("N73").ClearContents()
$ws.Range("H122").Value = 4446923
$ws.Range("I122").Value = 6062443
$ws.Range("J122").Value = 4244
$ws.Range("K122").Value = 18187329
$ws.Range("L122").Value = 12732
$ws.Range("M122").Value = -18184879
$ws.Range("N122").Value = -17632
$ws.Range("H132").Value = 5232.5
$ws.Range("I132").Value = 4240.7144
$ws.Range("J132").Value = 6224.2856
$ws.Range("K132").Value = 12722.1432
$ws.Range("L132").Value = 18672.8568
$ws.Range("M132").Value = -10192.1432
$ws.Range("N132").Value = -23732.8568
$ws.Range("H138").Value = 59332.832
$ws.Range("J138").Value = 59332.832
$ws.Range("L138").Value = 59332.832
$ws.Range("N138").Value = -69612.83199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 613.6087
$ws.Range("I16").Value = 622.0476
$ws.Range("K16").Value = 622.0476
$ws.Range("M16").Value = -452.0476
$ws.Range("H68").Value = 1725.5
$ws.Range("I68").Value = 2025
$ws.Range("J68").Value = 1597.1428
$ws.Range("K68").Value = 2025
$ws.Range("L68").Value = 1597.1428
$ws.Range("M68").Value = -1276
$ws.Range("N68").Value = -3095.1428
$ws.Range("H71").Value = 1725.5
$ws.Range("I71").Value = 2025
$ws.Range("J71").Value = 1597.1428
$ws.Range("K71").Value = 10125
$ws.Range("L71").Value = 7985.714
$ws.Range("M71").Value = -6381
$ws.Range("N71").Value = -15473.714
$ws.Range("H139").Value = 58682.668
$ws.Range("J139").Value = 58682.668
$ws.Range("L139").Value = 58682.668
$ws.Range("N139").Value = -68962.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 16213.444
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 16213.444
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 16213.444
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -16713.444
$ws.Range("H36").Value = 16213.444
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 16213.444
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 16213.444
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -16713.444
